$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (DNN vs cardio.)
$ws.Range("B2").Value = 0.002
$ws.Range("C2").Value = 0.739
$ws.Range("D2").Value = 0.083
$ws.Range("E2").Value = 0.414
$ws.Range("F2").Value = 0.527
$ws.Range("G2").Value = 0.366

# Row 3 (DNN vs emerg.)
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.285
$ws.Range("D3").Value = 0.025
$ws.Range("E3").Value = 0.257
$ws.Range("F3").Value = 0.763
$ws.Range("G3").Value = 0.739

# Row 4 (DNN vs stud.)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.025
$ws.Range("E4").Value = 0.058
$ws.Range("F4").Value = 0.617
$ws.Range("G4").Value = 0.285
